$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IPC")
$ws.Rows.Item(64).Delete()
